# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1. Update the "Date" metadata value to the new commit date.
# 2. On the "Elements" sheet, the two mapping columns (AK = "Mapping: RIM
#    Mapping", AL = "Mapping: Spécification métier vers l'extension ROR
#    AvailableTimeEffectiveOpeningClosingDate") are swapped: the contents
#    of column AK and column AL (header + every data row) are exchanged,
#    and the column widths follow the content (best-fit) so they are
#    swapped too.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Date" value (row 8, column B) ---------
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) ----------------
$elemWs = $wb.Worksheets.Item("Elements")

$lastRow = 10
$colAK = 37
$colAL = 38

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elemWs.Cells.Item($r, $colAK)
    $alCell = $elemWs.Cells.Item($r, $colAL)

    $akValue = $akCell.Value()
    $alValue = $alCell.Value()

    $akCell.Value = $alValue
    $alCell.Value = $akValue
}

# Column widths follow the (now swapped) content - the best-fit widths
# that used to belong to AL/AK are exchanged between the two columns too
# (AK becomes the wide "Spécification métier" column, AL becomes the
# narrow "RIM Mapping" column).
$elemWs.Columns.Item($colAK).ColumnWidth = 99
$elemWs.Columns.Item($colAL).ColumnWidth = 24.2
